$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lablist")

# --- Clear out the "T" column zero entries that were removed ---
$clearCells = @("T5","T6","T7","T8","T12","T15","T16","T17","T19","T20","T21","T22","T23","T35","T39","T40","T44","T45","T56","T57")
foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}

# --- Update / add "S" column values ---
$ws.Range("S31").Value = 1
$ws.Range("S33").Value = 0.75
$ws.Range("S50").Value = 1
$ws.Range("S51").Value = 0.75
$ws.Range("S56").Value = 0.5
$ws.Range("S58").Value = 0.5
$ws.Range("S60").Value = 1
$ws.Range("S67").Value = 1
$ws.Range("S70").Value = 1

# --- Update the active selection to match the saved view state ---
$ws.Activate()
$ws.Range("T57").Select()
